$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.018619775772095
$ws.Range("B1").Value = 6.237135887145996
$ws.Range("C1").Value = 3.238642692565918
$ws.Range("D1").Value = 1.402360320091248
$ws.Range("E1").Value = 0.9643964171409607
